$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-59: treasury delta data appended/re-derived from the original 15 rows.
# Column A for rows 16-43 stays numeric; rows 44-59 store the date code as text
# (matching the workbook authors original text-typed date column).
$rows = @(
    @{ Row = 16; A = "20071000"; B = 1642285; C = 1641080; D = -1205; AType = "num" }
    @{ Row = 17; A = "20071100"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 18; A = "20071200"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 19; A = "20071300"; B = 1641080; C = 1653555; D = 12475; AType = "num" }
    @{ Row = 20; A = "20071400"; B = 1653555; C = 1627709; D = -25846; AType = "num" }
    @{ Row = 21; A = "20071500"; B = 1627709; C = 1740455; D = 112746; AType = "num" }
    @{ Row = 22; A = "20071600"; B = 1740455; C = 1807305; D = 66850; AType = "num" }
    @{ Row = 23; A = "20071700"; B = 1807305; C = 1812801; D = 5496; AType = "num" }
    @{ Row = 24; A = "20071800"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 25; A = "20071900"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 26; A = "20072000"; B = 1812801; C = 1817270; D = 4469; AType = "num" }
    @{ Row = 27; A = "20072100"; B = 1817270; C = 1794448; D = -22822; AType = "num" }
    @{ Row = 28; A = "20072200"; B = 1794448; C = 1777180; D = -17268; AType = "num" }
    @{ Row = 29; A = "20070900"; B = 1624404; C = 1642285; D = 17881; AType = "num" }
    @{ Row = 30; A = "20071000"; B = 1642285; C = 1641080; D = -1205; AType = "num" }
    @{ Row = 31; A = "20071100"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 32; A = "20071200"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 33; A = "20071300"; B = 1641080; C = 1653555; D = 12475; AType = "num" }
    @{ Row = 34; A = "20071400"; B = 1653555; C = 1627709; D = -25846; AType = "num" }
    @{ Row = 35; A = "20071500"; B = 1627709; C = 1740455; D = 112746; AType = "num" }
    @{ Row = 36; A = "20071600"; B = 1740455; C = 1807305; D = 66850; AType = "num" }
    @{ Row = 37; A = "20071700"; B = 1807305; C = 1812801; D = 5496; AType = "num" }
    @{ Row = 38; A = "20071800"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 39; A = "20071900"; B = 0; C = 0; D = 0; AType = "num" }
    @{ Row = 40; A = "20072000"; B = 1812801; C = 1817270; D = 4469; AType = "num" }
    @{ Row = 41; A = "20072100"; B = 1817270; C = 1794448; D = -22822; AType = "num" }
    @{ Row = 42; A = "20072200"; B = 1794448; C = 1777180; D = -17268; AType = "num" }
    @{ Row = 43; A = "20072300"; B = 1777180; C = 1821824; D = 44644; AType = "num" }
    @{ Row = 44; A = "20070900"; B = 1624404; C = 1642285; D = 17881; AType = "text" }
    @{ Row = 45; A = "20071000"; B = 1642285; C = 1641080; D = -1205; AType = "text" }
    @{ Row = 46; A = "20071100"; B = 0; C = 0; D = 0; AType = "text" }
    @{ Row = 47; A = "20071200"; B = 0; C = 0; D = 0; AType = "text" }
    @{ Row = 48; A = "20071300"; B = 1641080; C = 1653555; D = 12475; AType = "text" }
    @{ Row = 49; A = "20071400"; B = 1653555; C = 1627709; D = -25846; AType = "text" }
    @{ Row = 50; A = "20071500"; B = 1627709; C = 1740455; D = 112746; AType = "text" }
    @{ Row = 51; A = "20071600"; B = 1740455; C = 1807305; D = 66850; AType = "text" }
    @{ Row = 52; A = "20071700"; B = 1807305; C = 1812801; D = 5496; AType = "text" }
    @{ Row = 53; A = "20071800"; B = 0; C = 0; D = 0; AType = "text" }
    @{ Row = 54; A = "20071900"; B = 0; C = 0; D = 0; AType = "text" }
    @{ Row = 55; A = "20072000"; B = 1812801; C = 1817270; D = 4469; AType = "text" }
    @{ Row = 56; A = "20072100"; B = 1817270; C = 1794448; D = -22822; AType = "text" }
    @{ Row = 57; A = "20072200"; B = 1794448; C = 1777180; D = -17268; AType = "text" }
    @{ Row = 58; A = "20072300"; B = 1777180; C = 1821824; D = 44644; AType = "text" }
    @{ Row = 59; A = "20072400"; B = 1821824; C = 1825498; D = 3674; AType = "text" }
)

foreach ($r in $rows) {
    if ($r.AType -eq "text") {
        $ws.Cells.Item($r.Row, 1).NumberFormat = "@"
        $ws.Cells.Item($r.Row, 1).Value = [string]$r.A
    } else {
        $ws.Cells.Item($r.Row, 1).Value = [double]$r.A
    }
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
